{"js": "// Remove the \"(letter and name)\" / \"(letters)\" qualifiers from the two\n// worksheet table header cells, per the commit message \"Updated lab 4\n// to remove letters\".\n//\n//   \"Most likely compound (letter and name)\" -> \"Most likely compound\"\n//   \"Other possibilities? (letters)\"         -> \"Other possibilities?\"\n\nconst replacements = [\n  [\"Most likely compound (letter and name)\", \"Most likely compound\"],\n  [\"Other possibilities? (letters)\", \"Other possibilities?\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const found = context.document.body.search(find, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"(letter and name)\" / \"(letters)\" qualifiers from the two\n# worksheet table header cells, per the commit message \"Updated lab 4\n# to remove letters\".\n#\n#   \"Most likely compound (letter and name)\" -> \"Most likely compound\"\n#   \"Other possibilities? (letters)\"         -> \"Other possibilities?\"\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"Most likely compound (letter and name)\", \"Most likely compound\"),\n    @(\"Other possibilities? (letters)\", \"Other possibilities?\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n"}
